$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run-boundary to persist at an absolute character position
# by toggling (and reverting) a character formatting attribute on the span
# [anchor, pos). This engine re-merges adjacent runs that share identical
# formatting unless a boundary has been "touched" by an edit, so growing the
# end of a range from a fixed, already-pinned anchor lets us pin a sequence
# of boundaries left-to-right without disturbing earlier ones.
# ---------------------------------------------------------------------------
function Pin-Boundary($anchor, $pos) {
    if ($pos -gt $anchor) {
        $r = $d.Range($anchor, $pos)
        $r.Font.Bold = 1
        $r.Font.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# Locate the "SSS-0036" paragraph by scanning the Paragraphs collection
# (robust to any header/footer paragraphs that may shift absolute indices).
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "SSS-0036*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the SSS-0036 paragraph"
}

$p36 = $d.Paragraphs.Item($targetIndex)
$paraStart = $p36.Range.Start
$paraEnd = $p36.Range.End

# ---------------------------------------------------------------------------
# Part 1: split the run
#   "– O sistema DEVE permitir que somente"
# into
#   "-"  and  " O sistema DEVE permitir que somente"
# while leaving the following pre-existing runs (" ", "o Administrador
# possa ", "alterar os dados de um Vendedor.") as separate runs.
# ---------------------------------------------------------------------------
$dashRange = $d.Range($paraStart, $paraEnd)
$dashFind = $dashRange.Find
$dashFound = $dashFind.Execute([string][char]0x2013, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $dashFound) {
    throw "Could not locate the en-dash in the SSS-0036 paragraph"
}
$dashStart = $dashRange.Start
$dashEnd = $dashRange.End

# Text following the dash, up to (not including) the next pre-existing run.
$afterDash = " O sistema DEVE permitir que somente"
$nextRun = "o Administrador possa "
$run1End = $dashEnd + $afterDash.Length
$spaceRunEnd = $run1End + 1
$thirdRunEnd = $spaceRunEnd + $nextRun.Length

# Replace the en-dash with a plain hyphen.
$dashRange.Text = "-"

# Pin: "-" | " O sistema DEVE permitir que somente" | " " | "o Administrador possa " | rest
Pin-Boundary $dashStart $dashEnd
Pin-Boundary $dashStart $run1End
Pin-Boundary $dashStart $spaceRunEnd
Pin-Boundary $dashStart $thirdRunEnd

# ---------------------------------------------------------------------------
# Part 2: insert a new paragraph after the SSS-0036 paragraph containing the
# new requirement SSS-0037, typed out as a sequence of small runs (mirrors
# how the text was actually authored).
# ---------------------------------------------------------------------------
$p36 = $d.Paragraphs.Item($targetIndex)
$p36.Range.InsertParagraphAfter()

$nextIndex = $targetIndex + 1
$p37 = $d.Paragraphs.Item($nextIndex)
$newParaStart = $p37.Range.Start

$fullText = "SSS-0037 - O sistema DEVE enviar a nota fiscal para o cliente."
$r = $d.Range($newParaStart, $newParaStart)
$r.InsertAfter($fullText)

$chunks = @('SSS-00', '3', '7', ' - O ', 'sistema', ' DEVE ', 'enviar', ' a nota fiscal', ' para o cliente', '.')
$pos = $newParaStart
foreach ($c in $chunks) {
    $pos = $pos + $c.Length
    Pin-Boundary $newParaStart $pos
}

Write-Host "Edit complete"
